$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62 (shifts existing rows 62..142 down to 63..143)
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new weekly price observation
$ws.Range("A62").Value = 4
$ws.Range("B62").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C62").Value = "Los Lagos"
$ws.Range("D62").Value = 44483
$ws.Range("E62").Value = 10
$ws.Range("F62").Value = 100112017
$ws.Range("G62").Value = "Apio"
$ws.Range("H62").Value = "Americana (o)"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 25
$ws.Range("K62").Value = 11000
$ws.Range("L62").Value = 11000
$ws.Range("M62").Value = 11000
$ws.Range("N62").Value = "`$/docena de matas"
$ws.Range("O62").Value = "Regi" + [char]0x00F3 + "n de Coquimbo"
$ws.Range("P62").Value = 1833
$ws.Range("Q62").Value = 6
$ws.Range("R62").Value = "Hortaliza"
